$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J190").Value = 'three carotenoid biosynthesis genes, chitobiose, D-glucoseaminide, fructose, glucose, mannose, fructan, rhamnulose, rhamnose, galactose, raffinose, stachyose, manninotriose, melibiose, xylose, sucrose, glucoside, cellulose, cellobiose, starch/glycogen, maltose'
$ws.Range("N190").Value = 'ammonia_assimilation, one nitric oxide reductase'
$ws.Range("P190").Value = 'partial sulfate_red_ass, thiosulfate'
$ws.Range("Q190").Value = 'xylose, iron, LPS export, lipoprotein release, phosphate, phospholipid/cholesterol, zinc'
$ws.Range("R190").Value = 'one chemotaxis protein'
$ws.Range("T190").Value = 'Oxidative phosphorylation'
$ws.Range("J191").Value = 'chitobiose, glucose, rhamnose, rhamnulose, lactose, galactan, raffinose, stachyose, manninotriose, melibiose, glucoside, cellobiose, starch/glycogen, maltose'
$ws.Range("N191").Value = 'ammonia_assimilation'
$ws.Range("P191").Value = 'tetrathionate, one sulfite reductase, extracellular sulfate transport'
$ws.Range("Q191").Value = 'xylose, iron, iron(III), LPS export, lipoprotein release, molybdate, oligopeptide, phosphate, phospholipid/cholesterol, sulfate, xylose'
$ws.Range("T191").Value = 'Oxidative phosphorylation'
$ws.Range("J192").Value = 'methanol, formate, chitobiose, fructose, glucose'
$ws.Range("N192").Value = 'nitrilase'
$ws.Range("P192").Value = 'partial sulfur oxidation (SOX)'
$ws.Range("Q192").Value = 'LPS export, lipoprotein release, phosphate, phospholipid, phospholipid/cholesterol'
$ws.Range("T192").Value = 'Oxidative phosphorylation'
$ws.Range("J193").Value = 'chitobiose, MurNAc, glucose, glycolate, 1,3-B-glucan, glucoside, cellobiose, starch/glycogen, maltose, trehalose'
$ws.Range("N193").Value = 'nitrogen fixation, ammonia_assimilation, hydroxylamine, one nitrite reductase'
$ws.Range("P193").Value = 'sulfate_red_dis, thiosulfate, trithioniate'
$ws.Range("Q193").Value = 'amino acid/amide, branched amino, cobalt/nickel, iron, LPS transport, microcin C, molybdate, phosphate, phospholipid, phospholipid/cholesterol, phosphonate, hydroxymethylpyrmidine, tungstate, zinc, type VI secretion'
$ws.Range("R193").Value = 'chemotaxis (purine?), flagellum'
$ws.Range("T193").Value = 'Oxidative phosphorylation'
$ws.Range("J194").Value = 'carotenoid biosynthesis, carbon fixation via Rubisco and reductive TCA, chitobiose, MurNAc, glucose, galacturonate, starch/glycogen, dextrin'
$ws.Range("N194").Value = 'nitrogen fixation, ammonia_assimilation, hydroxylamine'
$ws.Range("P194").Value = 'sulfate_red_ass, trithionate, one sulfite reductase from sulfate_dis_ass'
$ws.Range("Q194").Value = 'amino acid/amide, capsular polysaccharide, cobalt/nickel, iron, LPS export, lipoprotein release, molybdate, manganese/zinc/iron, macrolide, phosphate, phospholipid/cholesterol, sulfate'
$ws.Range("R194").Value = 'one chemotaxis protein'
$ws.Range("T194").Value = 'Oxidative phosphorylation'
$ws.Range("J195").Value = 'glucose, galactose, glycolate, glucoside, cellobiose, trehalose'
$ws.Range("N195").Value = 'ammonia_assimilation'
$ws.Range("P195").Value = 'taurine, sulfate_red_ass'
$ws.Range("Q195").Value = 'amino acid/amide, branched amino, carbohydrate, heme, iron, LPS transport, phosphate, sn-glycerol 3-phosphate, spermidine/putrescine'
$ws.Range("T195").Value = 'Oxidative phosphorylation'
$ws.Range("J196").Value = 'partial Wood-Ljungdahl, glucose, fructan, glucoside, cellobiose'
$ws.Range("N196").Value = 'ammonia_assimilation'
$ws.Range("Q196").Value = 'xylose, branched amino, iron (III), molybdate, monosaccharide, multiple sugar'
$ws.Range("T196").Value = 'Some oxidative phosphorylation'
$ws.Range("J197").Value = 'carbon fixation (RuBisCo), methanol, formate, acetate, MurNac, chitobiose, glycolate'
$ws.Range("N197").Value = 'ammonia_assimilation, nitrilase, nitronate monooxygenase'
$ws.Range("P197").Value = 'thiosulfate, partial sulfur oxidation (SOX)'
$ws.Range("Q197").Value = 'branched amino, heme, LPS export, lipoprotein release, phosphate, phospholipid/cholesterol, type IV secretion'
$ws.Range("R197").Value = 'two flagellum proteins'
$ws.Range("T197").Value = 'Oxidative phosphorylation'
$ws.Range("J198").Value = 'MurNAC, chitobiose, glycolate'
$ws.Range("N198").Value = 'ammonia_assimilation, nitrilase'
$ws.Range("Q198").Value = 'amino acid/amide, branched amino, heme, LPS export, lipoprotein release, microcin C, molybdate, paraquat-inducible, phosphate, phospholipid, phospholipid/cholesterol, tungstate'
$ws.Range("T198").Value = 'Some oxidative phosphorylation'
$ws.Range("J199").Value = 'chitobiose, glucose, galactose, rhamnulose, raffinose, stachyose, manninotriose, melibiose, glycolate, galacturonate, cellobiose, starch/glycogen, maltose, pectin, pectate'
$ws.Range("N199").Value = 'ammonia_assimilation, nitronate monooxygenase'
$ws.Range("P199").Value = 'one sulfite reductase, thiosulfate'
$ws.Range("Q199").Value = 'xylose, carbohydrate, iron, iron(III), LPS export, lipoprotein release, molybdate, oligopeptide, phosphate, phospholipid/cholesterol, sulfate, xylose'
$ws.Range("T199").Value = 'Oxidative phosphorylation'
$ws.Range("J200").Value = 'two carotenoid genes, glucose, glucoside, cellobiose, maltose'
$ws.Range("N200").Value = 'ammonia_assimilation'
$ws.Range("Q200").Value = 'branched amino, xylose, biotin, iron, iron (III), monosaccharide, multiple sugar, phosphate, ribose, teichoic acid, thiamine'
$ws.Range("T200").Value = 'Oxidative phosphorylation'
$ws.Range("J201").Value = 'methanol, formaldehyde, formate, MurNAC, chitobiose, glucose, galacturonate, 1,3-B-glucan, glucoside, cellobiose, starch/glycogen, trehalose, maltose'
$ws.Range("N201").Value = 'nitrilase'
$ws.Range("P201").Value = 'sulfate_red_ass, thiosulfate'
$ws.Range("Q201").Value = 'heme, LPS export, lipoprotein release, molybdate, molybdenum, phosphate, phospholipid/cholesterol, phosphonate, sulfate'
$ws.Range("R201").Value = 'chemotaxis (purine?), flagellum'
$ws.Range("T201").Value = 'Oxidative phosphorylation'
$ws.Range("J202").Value = 'methanol, formaldehyde, formate, MurNAc, chitobiose, glycolate, starch/glycogen'
$ws.Range("N202").Value = 'ammonia_assimilation'
$ws.Range("Q202").Value = 'LPS export, oligopeptide, phosphate, phospholipid, phospholipid/cholesterol, putrescine'
$ws.Range("T202").Value = 'Oxidative phosphorylation'
$ws.Range("J203").Value = 'partial Wood-Ljungdahl'
$ws.Range("L203").Value = 'chitobiose, glycolate, xylose, glucoside, cellobiose, cellulose'
$ws.Range("N203").Value = 'ammonia_assimilation'
$ws.Range("P203").Value = 'sulfur oxidation (SOX)'
$ws.Range("Q203").Value = 'carbohydrate, glycine betaine/proline, heme, lipoprotein release, molybdate, osmoprotectant, phosphate, phospholipid/cholesterol, phospholipid, putrescine, sodium, spermidine/putrescine'
$ws.Range("R203").Value = 'chemotaxis (purine?)'
$ws.Range("T203").Value = 'Oxidative phosphorylation'

$excel.ActiveWindow.ScrollRow = 159
$excel.ActiveWindow.ScrollColumn = 4
$ws.Range("A190:XFD203").Select()
